$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Gary Trent Jr.", "PG,SG", "Milwaukee Bucks"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Jared McCain", "PG,SG", "Philadelphia 76ers"),
    @("Harrison Barnes", "SF,PF", "San Antonio Spurs")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
